$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.093.09'
$ws.Range('D3').Value = '1.836.84'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6371'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.71%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07576'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.76%  '
$ws.Range('E9').Value = '  +1.42%  '
$ws.Range('E10').Value = '  +1.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07755'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.95%  '
$ws.Range('D12').Value = '1.843.41'
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.019'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6739'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '83.45'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009610'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.46%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.131'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.96%  '
$ws.Range('D18').Value = '29.134.17'
$ws.Range('E18').Value = '  +0.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.64'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '227.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.233'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '160.88'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.96%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1406'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.566'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.504'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.134'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.091'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.205'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05406'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.870'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7505'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.72%  '
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.663'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.57%  '
$ws.Range('D37').Value = '1.244.07'
$ws.Range('E37').Value = '  -2.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.763'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01793'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.628'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9074'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.65%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '102.46'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.95%  '
$ws.Range('D44').Value = '1.990.98'
$ws.Range('E44').Value = '  +0.73%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.28'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.28%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000123'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5114'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4116'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.096'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.798'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.57%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.652'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.05%  '
